$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing numeric cells on rows 8, 11, 12, 14 ---
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 7
$ws.Range("G8").Value = 0

$ws.Range("F11").Value = 8
$ws.Range("G11").Value = 0

$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 14
$ws.Range("G12").Value = 1

$ws.Range("F14").Value = 16
$ws.Range("G14").Value = 4
$ws.Range("E14").Value = 20

# --- Insert a new row at position 15, pushing the old row 15 (and the
#     hyperlink row further below) down by one ---
$ws.Rows("15").Insert()

# Populate the newly inserted row 15 with the "Presenter fgv-ek" / "Word" record
$ws.Range("A15").Value = "Presenter fgv-ek"
$ws.Range("B15").Value = "Word"
$ws.Range("D15").Value = 10
$ws.Range("F15").Value = 8
$ws.Range("H15").Value = "Pictori"

# --- Fix up the hyperlink that used to live on row 18 (now row 19) ---
$oldLink = $ws.Range("A18")
$oldLink.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A19"), "http://hungarian.joelonsoftware.com/Articles/PainlessSoftwareSchedules.html") | Out-Null

# --- Update the selected cell shown in the workbook ---
$ws.Range("J15").Select()
